# Auto-generated from the target diff: updates Price (D) and Volume(1h) (E)
# columns for the crypto listing rows on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.941.52'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '2.918.71'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'589.64"
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').Value = "'146.43"
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').Value = "'6.94"
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = "'33.57"
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = '3.402.10'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '60.963.91'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = "'6.69"
$ws.Range('E17').Value = '  -1.18%  '
$ws.Range('D18').Value = '2.920.74'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = "'432.39"
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('D20').Value = "'13.40"
$ws.Range('E20').Value = '  -1.55%  '
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('E22').Value = '  -0.34%  '
$ws.Range('D23').Value = "'81.30"
$ws.Range('E23').Value = '  +0.99%  '
$ws.Range('D24').Value = "'10.89"
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('D26').Value = "'11.85"
$ws.Range('E26').Value = '  -0.88%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  +4.53%  '
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').Value = "'6.97"
$ws.Range('E30').Value = '  -3.04%  '
$ws.Range('E31').Value = '  +3.06%  '
$ws.Range('D32').Value = "'26.64"
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = '0.0₃0863'
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('E35').Value = '  -0.40%  '
$ws.Range('E36').Value = '  -0.26%  '
$ws.Range('D37').Value = "'3.00"
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('E39').Value = '  -4.98%  '
$ws.Range('D40').Value = "'8.53"
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('D41').Value = "'41.52"
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('D42').Value = "'0.282"
$ws.Range('E42').Value = '  -4.68%  '
$ws.Range('D43').Value = "'376.30"
$ws.Range('E43').Value = '  -0.55%  '
$ws.Range('D44').Value = '2.706.61'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').Value = "'0.0343"
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('D46').Value = "'133.76"
$ws.Range('E46').Value = '  +0.88%  '
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('E48').Value = '  -4.20%  '
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('E50').Value = '  -3.09%  '
$ws.Range('E51').Value = '  -0.86%  '

# Re-apply the Normal style so the quote-prefix trick above does not
# leave these cells permanently marked as "text-quoted" (style-wise);
# only the underlying stored value needs to remain text.
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
